# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 52, pushing the existing
# rows 52 and 53 down to 53 and 54 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 52 (shifts old 52 -> 53, old 53 -> 54)
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new weekly entry
$ws.Range("A52").Value = 2
$ws.Range("B52").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C52").Value = "Coquimbo"
$ws.Range("D52").Value = 44595
$ws.Range("D52").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E52").Value = 4
$ws.Range("F52").Value = 100112032
$ws.Range("G52").Value = "Zapallo italiano"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 400
$ws.Range("K52").Value = 10000
$ws.Range("L52").Value = 12000
$ws.Range("M52").Value = 11000
$ws.Range("N52").Value = "$/caja 60 unidades"
$ws.Range("O52").Value = "Provincia de Limarí"
$ws.Range("P52").Value = 183
$ws.Range("Q52").Value = 60
$ws.Range("R52").Value = "Hortaliza"
